$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258; this shifts existing rows 258:316 down to 259:317
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with the latest week's record. Columns
# A,B,C,E,F,G,H,I,J,K,L,T carry the same market/product metadata as every
# other row in this sheet (Vega Modelo de Temuco / Frutilla); only D and
# M:S change with each new weekly observation.
$ws.Cells.Item(258, 1).Value = 10                         # A258 - Mercado ID
$ws.Cells.Item(258, 2).Value = "Vega Modelo de Temuco"    # B258 - Mercado
$ws.Cells.Item(258, 3).Value = "La Araucanía"             # C258 - Región
$ws.Cells.Item(258, 4).Value = 45034                      # D258 - Fecha
$ws.Cells.Item(258, 5).Value = 9                          # E258 - Codreg
$ws.Cells.Item(258, 6).Value = "Fruta"                    # F258 - Tipo
$ws.Cells.Item(258, 7).Value = 100101                     # G258 - Producto ID
$ws.Cells.Item(258, 8).Value = "Berries"                  # H258 - Producto
$ws.Cells.Item(258, 9).Value = 100112025                  # I258 - Categoría ID
$ws.Cells.Item(258, 10).Value = "Frutilla"                # J258 - Categoría
$ws.Cells.Item(258, 11).Value = "Sin especificar"         # K258 - Variedad
$ws.Cells.Item(258, 12).Value = "Primera"                 # L258 - Calidad
$ws.Cells.Item(258, 13).Value = 40                        # M258 - Volumen
$ws.Cells.Item(258, 14).Value = 8000                      # N258 - Precio mínimo
$ws.Cells.Item(258, 15).Value = 9000                      # O258 - Precio máximo
$ws.Cells.Item(258, 16).Value = 8500                      # P258 - Precio promedio ponderado
$ws.Cells.Item(258, 17).Value = "$/caja 7 kilos"          # Q258 - Unidad de comercialización
$ws.Cells.Item(258, 18).Value = "Región de La Araucanía"  # R258 - Origen
$ws.Cells.Item(258, 19).Value = 1214                      # S258 - Precio $/Kg
$ws.Cells.Item(258, 20).Value = 7                         # T258 - Kg / unidad

# Match the date-cell style used by every other row's "Fecha" column (D).
$ws.Cells.Item(258, 4).NumberFormat = $ws.Cells.Item(259, 4).NumberFormat
